$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.399.04'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.436.59'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.27'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.32'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +5.78%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.437.59'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.49%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.51%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.47%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.94%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.021.29'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.17'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.55%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.440.05'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.525.84'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.31'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +7.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.26'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '396.66'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +5.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.567'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.74'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.06%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.98%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.574.69'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.180'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.60'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.74%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.25'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.48%  '
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.18'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.09%  '
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.46'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -10.17%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.97'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.464.61'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.65%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.56'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.12'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '167.41'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0788'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.83'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.800'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.43%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.44%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '42.28'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.594.27'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.04%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.93'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.31%  '
